$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-12 from 45174 to 45175 (2023-09-06)
$ws.Range("C2:C12").Value = 45175
